# Update Leve profit figures (currentAveragePrice / HQ-NQ price & profit columns)
# across several item-source sheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 31304.666
$ws.Range("J108").Value = 31304.666
$ws.Range("L108").Value = 31304.666
$ws.Range("N108").Value = -38984.666

$ws.Range("H109").Value = 42680
$ws.Range("J109").Value = 42680
$ws.Range("L109").Value = 42680
$ws.Range("N109").Value = -45454

$ws.Range("H120").Value = 49714
$ws.Range("J120").Value = 49714
$ws.Range("L120").Value = 49714
$ws.Range("N120").Value = -59390

$ws.Range("H124").Value = 48517
$ws.Range("J124").Value = 48517
$ws.Range("L124").Value = 48517
$ws.Range("N124").Value = -58337

$ws.Range("H130").Value = 43160
$ws.Range("J130").Value = 43160
$ws.Range("L130").Value = 43160
$ws.Range("N130").Value = -53200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8856.679
$ws.Range("I32").Value = 7676.7183
$ws.Range("J32").Value = 15301.077
$ws.Range("K32").Value = 7676.7183
$ws.Range("L32").Value = 15301.077
$ws.Range("M32").Value = -7389.7183
$ws.Range("N32").Value = -15875.077

$ws.Range("H109").Value = 49092.25
$ws.Range("J109").Value = 49092.25
$ws.Range("L109").Value = 49092.25
$ws.Range("N109").Value = -51866.25

$ws.Range("H131").Value = 44285.668
$ws.Range("J131").Value = 44285.668
$ws.Range("L131").Value = 44285.668
$ws.Range("N131").Value = -54365.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 43332
$ws.Range("J116").Value = 43332
$ws.Range("L116").Value = 43332
$ws.Range("N116").Value = -52510

$ws.Range("H124").Value = 48966
$ws.Range("J124").Value = 48966
$ws.Range("L124").Value = 48966
$ws.Range("N124").Value = -58786

$ws.Range("H126").Value = 50780
$ws.Range("J126").Value = 50780
$ws.Range("L126").Value = 50780
$ws.Range("N126").Value = -60660

$ws.Range("H133").Value = 50997.75
$ws.Range("J133").Value = 50997.75
$ws.Range("L133").Value = 50997.75
$ws.Range("N133").Value = -61117.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 44718.145
$ws.Range("J20").Value = 44718.145
$ws.Range("L20").Value = 44718.145
$ws.Range("N20").Value = -45190.145

$ws.Range("H22").Value = 2610
$ws.Range("I22").Value = 360
$ws.Range("J22").Value = 4410
$ws.Range("K22").Value = 360
$ws.Range("L22").Value = 4410
$ws.Range("M22").Value = -10
$ws.Range("N22").Value = -5110

$ws.Range("H30").Value = 44718.145
$ws.Range("J30").Value = 44718.145
$ws.Range("L30").Value = 44718.145
$ws.Range("N30").Value = -44900.145

$ws.Range("H100").Value = 46931
$ws.Range("J100").Value = 46931
$ws.Range("L100").Value = 46931
$ws.Range("N100").Value = -49095

$ws.Range("H110").Value = 40563.332
$ws.Range("J110").Value = 40563.332
$ws.Range("L110").Value = 40563.332
$ws.Range("N110").Value = -48743.332

$ws.Range("H112").Value = 31895.25
$ws.Range("J112").Value = 31895.25
$ws.Range("L112").Value = 31895.25
$ws.Range("N112").Value = -34849.25

$ws.Range("H116").Value = 42874.668
$ws.Range("J116").Value = 42874.668
$ws.Range("L116").Value = 42874.668
$ws.Range("N116").Value = -52052.668

$ws.Range("H118").Value = 48742
$ws.Range("J118").Value = 48742
$ws.Range("L118").Value = 48742
$ws.Range("N118").Value = -52056

$ws.Range("H119").Value = 49376
$ws.Range("J119").Value = 49376
$ws.Range("L119").Value = 49376
$ws.Range("N119").Value = -59052

$ws.Range("H128").Value = 44718.145
$ws.Range("J128").Value = 44718.145
$ws.Range("L128").Value = 44718.145
$ws.Range("N128").Value = -54678.145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1338731.1
$ws.Range("I2").Value = 2409674
$ws.Range("J2").Value = 52.5
$ws.Range("K2").Value = 2409674
$ws.Range("L2").Value = 52.5
$ws.Range("M2").Value = -2409561
$ws.Range("N2").Value = -278.5

$ws.Range("H110").Value = 41781.332
$ws.Range("J110").Value = 41781.332
$ws.Range("L110").Value = 41781.332
$ws.Range("N110").Value = -49961.332

$ws.Range("H114").Value = 43025.25
$ws.Range("J114").Value = 43025.25
$ws.Range("L114").Value = 43025.25
$ws.Range("N114").Value = -51703.25

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H130").Value = 45474.855
$ws.Range("J130").Value = 45474.855
$ws.Range("L130").Value = 45474.855
$ws.Range("N130").Value = -55514.855

$ws.Range("H138").Value = 33155.91
$ws.Range("J138").Value = 33155.91
$ws.Range("L138").Value = 33155.91
$ws.Range("N138").Value = -43435.91

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 21773.273
$ws.Range("I2").Value = 1500
$ws.Range("J2").Value = 38667.668
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 38667.668
$ws.Range("M2").Value = -1388
$ws.Range("N2").Value = -38891.668

$ws.Range("H36").Value = 39538.4
$ws.Range("J36").Value = 39538.4
$ws.Range("L36").Value = 39538.4
$ws.Range("N36").Value = -40662.4

$ws.Range("H46").Value = 6100
$ws.Range("I46").Value = 4833.3335
$ws.Range("J46").Value = 8000
$ws.Range("K46").Value = 4833.3335
$ws.Range("L46").Value = 8000
$ws.Range("M46").Value = -4645.3335
$ws.Range("N46").Value = -8376

$ws.Range("H116").Value = 48680
$ws.Range("J116").Value = 48680
$ws.Range("L116").Value = 48680
$ws.Range("N116").Value = -57858

$ws.Range("H121").Value = 43416
$ws.Range("J121").Value = 43416
$ws.Range("L121").Value = 43416
$ws.Range("N121").Value = -46910

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 48636
$ws.Range("J110").Value = 48636
$ws.Range("L110").Value = 48636
$ws.Range("N110").Value = -56816

$ws.Range("H116").Value = 49680
$ws.Range("J116").Value = 49680
$ws.Range("L116").Value = 49680
$ws.Range("N116").Value = -58858

$ws.Range("H119").Value = 333363330
$ws.Range("J119").Value = 333363330
$ws.Range("L119").Value = 333363330
$ws.Range("N119").Value = -333373006

$ws.Range("H131").Value = 56490.668
$ws.Range("J131").Value = 56490.668
$ws.Range("L131").Value = 56490.668
$ws.Range("N131").Value = -66570.66800000001

